$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.668540666666667
$ws.Range("H2").Value = 14.005622
$ws.Range("I2").Value = 0.1221398332152395
$ws.Range("J2").Value = 0.1221398332152396
$ws.Range("M2").Value = 1.465773
$ws.Range("N2").Value = 4.397319
$ws.Range("O2").Value = 0.04255817041287053
$ws.Range("P2").Value = 0.04255817041287054
$ws.Range("Q2").Value = 6.843020858601999
$ws.Range("R2").Value = 61.587187727418
$ws.Range("S2").Value = 0.005198047836173748
$ws.Range("T2").Value = 0.00519804783617375
$ws.Range("G3").Value = 4.668540666666667
$ws.Range("H3").Value = 14.005622
$ws.Range("I3").Value = 0.1221398332152395
$ws.Range("J3").Value = 0.1221398332152396
$ws.Range("O3").Value = 0.2518776232723828
$ws.Range("P3").Value = 0.2518776232723829
$ws.Range("Q3").Value = 40.49995131714488
$ws.Range("R3").Value = 364.499561854304
$ws.Range("S3").Value = 0.03076429089713978
$ws.Range("T3").Value = 0.03076429089713979
$ws.Range("G4").Value = 4.668540666666667
$ws.Range("H4").Value = 14.005622
$ws.Range("I4").Value = 0.1221398332152395
$ws.Range("J4").Value = 0.1221398332152396
$ws.Range("M4").Value = 3.352351333333333
$ws.Range("N4").Value = 10.057054
$ws.Range("O4").Value = 0.09733426617069202
$ws.Range("P4").Value = 0.09733426617069205
$ws.Range("Q4").Value = 15.65058852862089
$ws.Range("R4").Value = 140.855296757588
$ws.Range("S4").Value = 0.01188839103621605
$ws.Range("T4").Value = 0.01188839103621606
$ws.Range("G5").Value = 4.668540666666667
$ws.Range("H5").Value = 14.005622
$ws.Range("I5").Value = 0.1221398332152395
$ws.Range("J5").Value = 0.1221398332152396
$ws.Range("M5").Value = 10.935394
$ws.Range("N5").Value = 32.806182
$ws.Range("O5").Value = 0.3175050716474392
$ws.Range("P5").Value = 0.3175050716474393
$ws.Range("Q5").Value = 51.05233159502267
$ws.Range("R5").Value = 459.470984355204
$ws.Range("S5").Value = 0.03878001649601091
$ws.Range("T5").Value = 0.03878001649601092
$ws.Range("G6").Value = 4.668540666666667
$ws.Range("H6").Value = 14.005622
$ws.Range("I6").Value = 0.1221398332152395
$ws.Range("J6").Value = 0.1221398332152396
$ws.Range("M6").Value = 10.01304
$ws.Range("N6").Value = 30.03912
$ws.Range("O6").Value = 0.2907248684966152
$ws.Range("P6").Value = 0.2907248684966153
$ws.Range("Q6").Value = 46.74628443696
$ws.Range("R6").Value = 420.7165599326401
$ws.Range("S6").Value = 0.03550908694969903
$ws.Range("T6").Value = 0.03550908694969904
$ws.Range("I7").Value = 0.6401930168419905
$ws.Range("J7").Value = 0.6401930168419905
$ws.Range("M7").Value = 1.465773
$ws.Range("N7").Value = 4.397319
$ws.Range("O7").Value = 0.04255817041287053
$ws.Range("P7").Value = 0.04255817041287054
$ws.Range("Q7").Value = 35.86753029260299
$ws.Range("R7").Value = 322.807772633427
$ws.Range("S7").Value = 0.02724544350789112
$ws.Range("T7").Value = 0.02724544350789113
$ws.Range("I8").Value = 0.6401930168419905
$ws.Range("J8").Value = 0.6401930168419905
$ws.Range("O8").Value = 0.2518776232723828
$ws.Range("P8").Value = 0.2518776232723829
$ws.Range("S8").Value = 0.1612502955177371
$ws.Range("T8").Value = 0.1612502955177371
$ws.Range("I9").Value = 0.6401930168419905
$ws.Range("J9").Value = 0.6401930168419905
$ws.Range("M9").Value = 3.352351333333333
$ws.Range("N9").Value = 10.057054
$ws.Range("O9").Value = 0.09733426617069202
$ws.Range("P9").Value = 0.09733426617069205
$ws.Range("Q9").Value = 82.03218574757578
$ws.Range("R9").Value = 738.289671728182
$ws.Range("S9").Value = 0.06231271750191662
$ws.Range("T9").Value = 0.06231271750191664
$ws.Range("I10").Value = 0.6401930168419905
$ws.Range("J10").Value = 0.6401930168419905
$ws.Range("M10").Value = 10.935394
$ws.Range("N10").Value = 32.806182
$ws.Range("O10").Value = 0.3175050716474392
$ws.Range("P10").Value = 0.3175050716474393
$ws.Range("Q10").Value = 267.5895759824674
$ws.Range("R10").Value = 2408.306183842206
$ws.Range("S10").Value = 0.2032645296806065
$ws.Range("T10").Value = 0.2032645296806065
$ws.Range("I11").Value = 0.6401930168419905
$ws.Range("J11").Value = 0.6401930168419905
$ws.Range("M11").Value = 10.01304
$ws.Range("N11").Value = 30.03912
$ws.Range("O11").Value = 0.2907248684966152
$ws.Range("P11").Value = 0.2907248684966153
$ws.Range("Q11").Value = 245.01953271144
$ws.Range("R11").Value = 2205.17579440296
$ws.Range("S11").Value = 0.1861200306338391
$ws.Range("T11").Value = 0.1861200306338391
$ws.Range("G12").Value = 3.558094333333333
$ws.Range("H12").Value = 10.674283
$ws.Range("I12").Value = 0.09308798604676512
$ws.Range("J12").Value = 0.09308798604676513
$ws.Range("M12").Value = 1.465773
$ws.Range("N12").Value = 4.397319
$ws.Range("O12").Value = 0.04255817041287053
$ws.Range("P12").Value = 0.04255817041287054
$ws.Range("Q12").Value = 5.215358605252999
$ws.Range("R12").Value = 46.93822744727699
$ws.Range("S12").Value = 0.003961654373569144
$ws.Range("T12").Value = 0.003961654373569145
$ws.Range("G13").Value = 3.558094333333333
$ws.Range("H13").Value = 10.674283
$ws.Range("I13").Value = 0.09308798604676512
$ws.Range("J13").Value = 0.09308798604676513
$ws.Range("O13").Value = 0.2518776232723828
$ws.Range("P13").Value = 0.2518776232723829
$ws.Range("Q13").Value = 30.86674350096177
$ws.Range("R13").Value = 277.800691508656
$ws.Range("S13").Value = 0.02344678068067193
$ws.Range("T13").Value = 0.02344678068067194
$ws.Range("G14").Value = 3.558094333333333
$ws.Range("H14").Value = 10.674283
$ws.Range("I14").Value = 0.09308798604676512
$ws.Range("J14").Value = 0.09308798604676513
$ws.Range("M14").Value = 3.352351333333333
$ws.Range("N14").Value = 10.057054
$ws.Range("O14").Value = 0.09733426617069202
$ws.Range("P14").Value = 0.09733426617069205
$ws.Range("Q14").Value = 11.92798228247578
$ws.Range("R14").Value = 107.351840542282
$ws.Range("S14").Value = 0.009060650811169501
$ws.Range("T14").Value = 0.009060650811169505
$ws.Range("G15").Value = 3.558094333333333
$ws.Range("H15").Value = 10.674283
$ws.Range("I15").Value = 0.09308798604676512
$ws.Range("J15").Value = 0.09308798604676513
$ws.Range("M15").Value = 10.935394
$ws.Range("N15").Value = 32.806182
$ws.Range("O15").Value = 0.3175050716474392
$ws.Range("P15").Value = 0.3175050716474393
$ws.Range("Q15").Value = 38.90916342416733
$ws.Range("R15").Value = 350.182470817506
$ws.Range("S15").Value = 0.02955590767929398
$ws.Range("T15").Value = 0.02955590767929399
$ws.Range("G16").Value = 3.558094333333333
$ws.Range("H16").Value = 10.674283
$ws.Range("I16").Value = 0.09308798604676512
$ws.Range("J16").Value = 0.09308798604676513
$ws.Range("M16").Value = 10.01304
$ws.Range("N16").Value = 30.03912
$ws.Range("O16").Value = 0.2907248684966152
$ws.Range("P16").Value = 0.2907248684966153
$ws.Range("Q16").Value = 35.62734088344
$ws.Range("R16").Value = 320.64606795096
$ws.Range("S16").Value = 0.02706299250206054
$ws.Range("T16").Value = 0.02706299250206055
$ws.Range("G17").Value = 1.877690666666667
$ws.Range("H17").Value = 5.633072
$ws.Range("I17").Value = 0.04912473537908105
$ws.Range("J17").Value = 0.04912473537908105
$ws.Range("M17").Value = 1.465773
$ws.Range("N17").Value = 4.397319
$ws.Range("O17").Value = 0.04255817041287053
$ws.Range("P17").Value = 0.04255817041287054
$ws.Range("Q17").Value = 2.752268281552
$ws.Range("R17").Value = 24.770414533968
$ws.Range("S17").Value = 0.002090658859750101
$ws.Range("T17").Value = 0.002090658859750101
$ws.Range("G18").Value = 1.877690666666667
$ws.Range("H18").Value = 5.633072
$ws.Range("I18").Value = 0.04912473537908105
$ws.Range("J18").Value = 0.04912473537908105
$ws.Range("O18").Value = 0.2518776232723828
$ws.Range("P18").Value = 0.2518776232723829
$ws.Range("Q18").Value = 16.28911174141156
$ws.Range("R18").Value = 146.602005672704
$ws.Range("S18").Value = 0.01237342159116767
$ws.Range("T18").Value = 0.01237342159116767
$ws.Range("G19").Value = 1.877690666666667
$ws.Range("H19").Value = 5.633072
$ws.Range("I19").Value = 0.04912473537908105
$ws.Range("J19").Value = 0.04912473537908105
$ws.Range("M19").Value = 3.352351333333333
$ws.Range("N19").Value = 10.057054
$ws.Range("O19").Value = 0.09733426617069202
$ws.Range("P19").Value = 0.09733426617069205
$ws.Range("Q19").Value = 6.294678809987556
$ws.Range("R19").Value = 56.65210928988801
$ws.Range("S19").Value = 0.004781520068952286
$ws.Range("T19").Value = 0.004781520068952287
$ws.Range("G20").Value = 1.877690666666667
$ws.Range("H20").Value = 5.633072
$ws.Range("I20").Value = 0.04912473537908105
$ws.Range("J20").Value = 0.04912473537908105
$ws.Range("M20").Value = 10.935394
$ws.Range("N20").Value = 32.806182
$ws.Range("O20").Value = 0.3175050716474392
$ws.Range("P20").Value = 0.3175050716474393
$ws.Range("Q20").Value = 20.53328725012267
$ws.Range("R20").Value = 184.799585251104
$ws.Range("S20").Value = 0.01559735262619662
$ws.Range("T20").Value = 0.01559735262619662
$ws.Range("G21").Value = 1.877690666666667
$ws.Range("H21").Value = 5.633072
$ws.Range("I21").Value = 0.04912473537908105
$ws.Range("J21").Value = 0.04912473537908105
$ws.Range("M21").Value = 10.01304
$ws.Range("N21").Value = 30.03912
$ws.Range("O21").Value = 0.2907248684966152
$ws.Range("P21").Value = 0.2907248684966153
$ws.Range("Q21").Value = 18.80139175296
$ws.Range("R21").Value = 169.21252577664
$ws.Range("S21").Value = 0.01428178223301436
$ws.Range("T21").Value = 0.01428178223301436
$ws.Range("G22").Value = 3.648546666666667
$ws.Range("H22").Value = 10.94564
$ws.Range("I22").Value = 0.09545442851692375
$ws.Range("J22").Value = 0.09545442851692375
$ws.Range("M22").Value = 1.465773
$ws.Range("N22").Value = 4.397319
$ws.Range("O22").Value = 0.04255817041287053
$ws.Range("P22").Value = 0.04255817041287054
$ws.Range("Q22").Value = 5.34794119324
$ws.Range("R22").Value = 48.13147073916
$ws.Range("S22").Value = 0.004062365835486409
$ws.Range("T22").Value = 0.00406236583548641
$ws.Range("G23").Value = 3.648546666666667
$ws.Range("H23").Value = 10.94564
$ws.Range("I23").Value = 0.09545442851692375
$ws.Range("J23").Value = 0.09545442851692375
$ws.Range("O23").Value = 0.2518776232723828
$ws.Range("P23").Value = 0.2518776232723829
$ws.Range("Q23").Value = 31.65142448760889
$ws.Range("R23").Value = 284.86282038848
$ws.Range("S23").Value = 0.02404283458566632
$ws.Range("T23").Value = 0.02404283458566632
$ws.Range("G24").Value = 3.648546666666667
$ws.Range("H24").Value = 10.94564
$ws.Range("I24").Value = 0.09545442851692375
$ws.Range("J24").Value = 0.09545442851692375
$ws.Range("M24").Value = 3.352351333333333
$ws.Range("N24").Value = 10.057054
$ws.Range("O24").Value = 0.09733426617069202
$ws.Range("P24").Value = 0.09733426617069205
$ws.Range("Q24").Value = 12.23121028272889
$ws.Range("R24").Value = 110.08089254456
$ws.Range("S24").Value = 0.009290986752437551
$ws.Range("T24").Value = 0.009290986752437553
$ws.Range("G25").Value = 3.648546666666667
$ws.Range("H25").Value = 10.94564
$ws.Range("I25").Value = 0.09545442851692375
$ws.Range("J25").Value = 0.09545442851692375
$ws.Range("M25").Value = 10.935394
$ws.Range("N25").Value = 32.806182
$ws.Range("O25").Value = 0.3175050716474392
$ws.Range("P25").Value = 0.3175050716474393
$ws.Range("Q25").Value = 39.89829532738667
$ws.Range("R25").Value = 359.08465794648
$ws.Range("S25").Value = 0.03030726516533124
$ws.Range("T25").Value = 0.03030726516533125
$ws.Range("G26").Value = 3.648546666666667
$ws.Range("H26").Value = 10.94564
$ws.Range("I26").Value = 0.09545442851692375
$ws.Range("J26").Value = 0.09545442851692375
$ws.Range("M26").Value = 10.01304
$ws.Range("N26").Value = 30.03912
$ws.Range("O26").Value = 0.2907248684966152
$ws.Range("P26").Value = 0.2907248684966153
$ws.Range("Q26").Value = 46.74628443696
$ws.Range("R26").Value = 328.7973934368
$ws.Range("S26").Value = 0.02775097617800221
$ws.Range("T26").Value = 0.02775097617800222
